$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G5").Value = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G6").Value = "mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G7").Value = "menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G11").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G13").Value = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G20").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G24").Value = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G25").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
